$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.784.08'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.088.42'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.29'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.85%  '

$ws.Range("E12").Value = '  +1.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.396.86'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.20'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.778'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.36'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.087.26'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.737.11'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.12'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.71'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.70'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.46'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.46%  '

$ws.Range("E28").Value = '  -3.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.54'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.19%  '

$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0638'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.61'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.82'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("E37").Value = '  -1.78%  '

$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.37'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("E40").Value = '  +10.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.44'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0969'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.88%  '

$ws.Range("E44").Value = '  +1.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.68'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.452.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.11'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.57%  '

$ws.Range("E48").Value = '  -0.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.22'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.44%  '

$ws.Range("E50").Value = '  -1.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.280.83'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.37%  '
